$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 2
$ws.Range("E2").Value = 93

# Row 10
$ws.Range("E10").Value = 513
$ws.Range("F10").Value = 258
$ws.Range("H10").Value = 353

# Row 11
$ws.Range("E11").Value = 334
$ws.Range("G11").Value = 67
$ws.Range("H11").Value = 253

# Row 12
$ws.Range("E12").Value = 501
$ws.Range("G12").Value = 84
$ws.Range("H12").Value = 357

# Row 13
$ws.Range("F13").Value = 69
$ws.Range("H13").Value = 103

# Row 23
$ws.Range("F23").Value = 93
$ws.Range("H23").Value = 144

# Row 27
$ws.Range("F27").Value = 159
$ws.Range("H27").Value = 240

# Row 28
$ws.Range("E28").Value = 192

# Row 30
$ws.Range("E30").Value = 200
$ws.Range("F30").Value = 120
$ws.Range("H30").Value = 172

# Row 33
$ws.Range("E33").Value = 284

# Row 34
$ws.Range("F34").Value = 140
$ws.Range("H34").Value = 179

# Row 35
$ws.Range("F35").Value = 91
$ws.Range("H35").Value = 118

# Row 36
$ws.Range("E36").Value = 69

# Row 40
$ws.Range("F40").Value = 121
$ws.Range("H40").Value = 201

# Row 47
$ws.Range("F47").Value = 222
$ws.Range("H47").Value = 314

# Row 50
$ws.Range("E50").Value = 237
$ws.Range("F50").Value = 112
$ws.Range("H50").Value = 185
